# Updated symbol list (crypto price/volume refresh) -- applies the cell-level
# text updates described by the commit diff. Every touched cell in the sheet
# stores its value as literal text (e.g. "310.64", "1.56%"), so we force the
# cell's number format to Text ("@") before writing the new value. Otherwise
# Excel's normal autodetection would coerce a numeric-looking string into a
# real number, or a "NN.NN%" string into a percentage value -- changing the
# cell's type/content instead of just refreshing its text, which is not what
# the source diff shows (every changed cell stays an inline/text string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($sheet, $addr, $text)
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

Set-CellText $ws "D2" "310.64"
Set-CellText $ws "E2" "1.56%"
Set-CellText $ws "E3" "-1.88%"
Set-CellText $ws "E4" "0.84%"
Set-CellText $ws "D5" "0.08231"
Set-CellText $ws "E5" "4.88%"
Set-CellText $ws "D6" "2.068"
Set-CellText $ws "E6" "-8.26%"
Set-CellText $ws "D7" "7.946"
Set-CellText $ws "E7" "-0.36%"
Set-CellText $ws "E8" "11.33%"
Set-CellText $ws "E9" "-0.20%"
Set-CellText $ws "D10" "0.1091"
Set-CellText $ws "E10" "12.04%"
Set-CellText $ws "D11" "0.1918"
Set-CellText $ws "E11" "3.03%"
Set-CellText $ws "D12" "0.09404"
Set-CellText $ws "E12" "5.26%"
Set-CellText $ws "D13" "0.03641"
Set-CellText $ws "E13" "-3.74%"
Set-CellText $ws "D14" "0.09919"
Set-CellText $ws "E14" "0.13%"
Set-CellText $ws "D15" "0.001435"
Set-CellText $ws "E15" "-0.87%"
Set-CellText $ws "D16" "0.005862"
Set-CellText $ws "E16" "3.26%"
Set-CellText $ws "E17" "0.01%"
Set-CellText $ws "D18" "4.126"
Set-CellText $ws "E18" "-0.71%"
Set-CellText $ws "D19" "0.3429"
Set-CellText $ws "E19" "0.22%"
Set-CellText $ws "E20" "-1.04%"
Set-CellText $ws "D21" "5.088"
Set-CellText $ws "E21" "-1.31%"
Set-CellText $ws "D22" "0.2205"
Set-CellText $ws "E22" "-2.39%"
Set-CellText $ws "D23" "0.04545"
Set-CellText $ws "E23" "-0.70%"
Set-CellText $ws "E24" "-0.82%"
Set-CellText $ws "D25" "0.004815"
Set-CellText $ws "E25" "1.00%"
Set-CellText $ws "D26" "0.0001249"
Set-CellText $ws "E26" "-4.39%"
Set-CellText $ws "D27" "0.0004448"
Set-CellText $ws "E27" "-6.12%"
Set-CellText $ws "D39" "0.01986"
Set-CellText $ws "E39" "2.25%"
Set-CellText $ws "D40" "0.04934"
Set-CellText $ws "E40" "-0.66%"
Set-CellText $ws "D41" "0.007704"
Set-CellText $ws "E41" "-1.25%"
Set-CellText $ws "D42" "0.009977"
Set-CellText $ws "E42" "27.29%"
Set-CellText $ws "D43" "0.1386"
Set-CellText $ws "E43" "-0.11%"
Set-CellText $ws "D44" "0.002113"
Set-CellText $ws "E44" "-1.72%"
Set-CellText $ws "E45" "2.86%"
Set-CellText $ws "D46" "0.00006563"
Set-CellText $ws "E46" "6.26%"
Set-CellText $ws "E47" "-0.48%"
Set-CellText $ws "D48" "61.86"
Set-CellText $ws "E48" "19.52%"
Set-CellText $ws "E49" "-21.49%"
Set-CellText $ws "E50" "-0.48%"
Set-CellText $ws "E51" "-0.48%"
